$d = $word.ActiveDocument

function Add-PlainRun($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $r = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $r.InsertAfter($text)
    $r.Font.Name = "Californian FB"
}

function Add-BoldRun($paraIndex, $text) {
    $p = $d.Paragraphs($paraIndex)
    $r = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $r.InsertAfter($text)
    $r.Font.Name = "Californian FB"
    $r.Font.Bold = 1
    $r.Font.BoldBi = 1
}

# Locate the paragraph ending with "This is the ascending limb. " and
# append: "Its walls are " + "impermeable" (bold) + " to water."
$ascendingParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "This is the ascending limb. `r") {
        $ascendingParaIndex = $i
        break
    }
}

Add-PlainRun $ascendingParaIndex "Its walls are "
Add-BoldRun $ascendingParaIndex "impermeable"
Add-PlainRun $ascendingParaIndex " to water."

# Locate the paragraph ending with "This is the descending limb." and
# append: " Its walls are " + "permeable" (bold) + " to water."
$descendingParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "This is the descending limb.`r") {
        $descendingParaIndex = $i
        break
    }
}

Add-PlainRun $descendingParaIndex " Its walls are "
Add-BoldRun $descendingParaIndex "permeable"
Add-PlainRun $descendingParaIndex " to water."

Write-Host "Ascending limb paragraph index:" $ascendingParaIndex
Write-Host "Descending limb paragraph index:" $descendingParaIndex
Write-Host "Ascending text now:" $d.Paragraphs($ascendingParaIndex).Range.Text
Write-Host "Descending text now:" $d.Paragraphs($descendingParaIndex).Range.Text
